$d = $word.ActiveDocument

function Find-AllOccurrences($text, $needle) {
    $positions = @()
    $idx = 0
    while ($true) {
        $found = $text.IndexOf($needle, $idx)
        if ($found -lt 0) { break }
        $positions += $found
        $idx = $found + 1
    }
    return $positions
}

# -----------------------------------------------------------------
# Step 1: fix the "Track4Help" -> "Data4Help" typo (3 occurrences).
# Word's "type over a selection" behaviour leaves the replaced text in
# its own run while the untouched remainder ("4Help...") stays merged
# with whatever follows it. We reproduce that precisely by temporarily
# dropping bookmarks around "Track" to fence off the run-merge that
# this engine performs whenever a Range.Text edit touches a run.
# -----------------------------------------------------------------
$full = $d.Content.Text
$positions = Find-AllOccurrences $full "Track4Help"

for ($i = $positions.Count - 1; $i -ge 0; $i--) {
    $p = $positions[$i]

    $leftBm = "zzFence" + $i + "L"
    $rightBm = "zzFence" + $i + "R"

    $rLeft = $d.Range($p, $p)
    $d.Bookmarks.Add($leftBm, $rLeft)
    $rRight = $d.Range($p + 5, $p + 5)
    $d.Bookmarks.Add($rightBm, $rRight)

    $r = $d.Range($p, $p + 5)
    $r.Text = "Data"

    $d.Bookmarks($leftBm).Delete()
    $d.Bookmarks($rightBm).Delete()
}

# -----------------------------------------------------------------
# Step 2: the _GoBack bookmark (Word's "last edit location" marker)
# moves to the newest edit point: right after "software" in
# "...system is designed as a software application...".
# -----------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# -----------------------------------------------------------------
# Step 3: with the old _GoBack bookmark gone, the two runs it used to
# separate ("...retrieve the d" / "ata of the device once a day.")
# become plain adjacent text and should merge back into a single run,
# the same way Word recombines a run it had split purely to host the
# now-removed bookmark. We again fence the merge so it only touches
# those two runs and does not bleed into the neighbours ("The " and
# " In order to know...").
# -----------------------------------------------------------------
$full3 = $d.Content.Text
$mergeStartText = "second service is meant to call help for elderly people if they need it. Thus the application should monitor the data continuously, and not just retrieve the data of the device once a day."
$mergeStart = $full3.IndexOf("second service is meant")
$mergeEnd = $mergeStart + $mergeStartText.Length

$rFenceL = $d.Range($mergeStart, $mergeStart)
$d.Bookmarks.Add("zzMergeFenceL", $rFenceL)
$rFenceR = $d.Range($mergeEnd, $mergeEnd)
$d.Bookmarks.Add("zzMergeFenceR", $rFenceR)

$idxD = $full3.IndexOf("retrieve the d")
$posD = $idxD + "retrieve the ".Length
$rTouch = $d.Range($posD, $posD + 1)
$rTouch.Text = "e"
$rTouch2 = $d.Range($posD, $posD + 1)
$rTouch2.Text = "d"

$d.Bookmarks("zzMergeFenceL").Delete()
$d.Bookmarks("zzMergeFenceR").Delete()

# -----------------------------------------------------------------
# Step 4: re-add _GoBack at the new edit location (after "software").
# -----------------------------------------------------------------
$full4 = $d.Content.Text
$idx4 = $full4.IndexOf("as a software application")
$newPos = $idx4 + "as a software".Length
$newRange = $d.Range($newPos, $newPos)
$d.Bookmarks.Add("_GoBack", $newRange)
